$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0)
$ws.Range("B3").Value = 0.2480681570822581
$ws.Range("C3").Value = 0.9558909634454847
$ws.Range("D3").Value = 1.537345962423154
$ws.Range("E3").Value = 1.239897561261879
$ws.Range("F3").Value = 1.242131370099838
$ws.Range("G3").Value = 23

# Row 4 (Q1)
$ws.Range("B4").Value = 0.2672114150504847
$ws.Range("C4").Value = 1.59686641897852
$ws.Range("D4").Value = 10.20774690982041
$ws.Range("E4").Value = 3.194956480113682
$ws.Range("F4").Value = 3.258685020567687
$ws.Range("G4").Value = 22
